$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the edited range to Text format first so numeric-looking strings
# (e.g. "26.738.36", "4.23") are stored as literal text, matching the
# source workbooks inline-string cells instead of being parsed as numbers.
$editRange = $ws.Range("B2:E51")
$editRange.NumberFormat = "@"

$ws.Range('D2').Value = '26.738.36'
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').Value = '1.649.32'
$ws.Range('E3').Value = '  +0.75%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '215.87'
$ws.Range('E5').Value = '  +1.21%  '
$ws.Range('E6').Value = '  +1.17%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '0.254'
$ws.Range('E8').Value = '  -0.19%  '
$ws.Range('D9').Value = '0.0630'
$ws.Range('E9').Value = '  +1.01%  '
$ws.Range('D10').Value = '19.43'
$ws.Range('E10').Value = '  +1.37%  '
$ws.Range('D11').Value = '0.0845'
$ws.Range('E11').Value = '  +0.46%  '
$ws.Range('D12').Value = '1.880.68'
$ws.Range('E12').Value = '  +0.84%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.674.87'
$ws.Range('E13').Value = '  +2.19%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '4.23'
$ws.Range('E14').Value = '  +3.11%  '
$ws.Range('D15').Value = '0.536'
$ws.Range('E15').Value = '  +1.47%  '
$ws.Range('D16').Value = '66.68'
$ws.Range('E16').Value = '  +5.28%  '
$ws.Range('D17').Value = '26.809.12'
$ws.Range('E17').Value = '  +0.54%  '
$ws.Range('E18').Value = '  +1.76%  '
$ws.Range('D19').Value = '221.50'
$ws.Range('E19').Value = '  +1.32%  '
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('E21').Value = '  +2.45%  '
$ws.Range('D22').Value = '6.38'
$ws.Range('E22').Value = '  +2.57%  '
$ws.Range('E23').Value = '  +1.07%  '
$ws.Range('D24').Value = '2.16'
$ws.Range('E24').Value = '  +12.18%  '
$ws.Range('D25').Value = '147.46'
$ws.Range('E25').Value = '  -1.10%  '
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('E27').Value = '  +0.45%  '
$ws.Range('D28').Value = '7.13'
$ws.Range('E28').Value = '  +4.08%  '
$ws.Range('D29').Value = '15.97'
$ws.Range('E29').Value = '  +3.44%  '
$ws.Range('D30').Value = '0.0524'
$ws.Range('E30').Value = '  +1.67%  '
$ws.Range('D31').Value = '1.18'
$ws.Range('E31').Value = '  +0.69%  '
$ws.Range('D32').Value = '3.44'
$ws.Range('E32').Value = '  +3.91%  '
$ws.Range('D33').Value = '3.07'
$ws.Range('E33').Value = '  +4.45%  '
$ws.Range('D34').Value = '1.57'
$ws.Range('E34').Value = '  +3.98%  '
$ws.Range('D35').Value = '1.295.60'
$ws.Range('E35').Value = '  +8.03%  '
$ws.Range('D36').Value = '0.0184'
$ws.Range('E36').Value = '  +6.22%  '
$ws.Range('E37').Value = '  +0.97%  '
$ws.Range('D38').Value = '0.833'
$ws.Range('E38').Value = '  +3.02%  '
$ws.Range('E39').Value = '  +4.52%  '
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').Value = '0.816'
$ws.Range('E41').Value = '  +2.74%  '
$ws.Range('E43').Value = '  +0.92%  '
$ws.Range('D44').Value = '1.791.98'
$ws.Range('D45').Value = '93.92'
$ws.Range('E45').Value = '  +1.79%  '
$ws.Range('D46').Value = '60.24'
$ws.Range('E46').Value = '  +9.77%  '
$ws.Range('E47').Value = '  +5.74%  '
$ws.Range('D48').Value = '0.0518'
$ws.Range('E48').Value = '  +1.08%  '
$ws.Range('D49').Value = '7.83'
$ws.Range('E49').Value = '  +2.63%  '
$ws.Range('D50').Value = '0.0983'
$ws.Range('E50').Value = '  +3.80%  '
$ws.Range('E51').Value = '  -0.72%  '

# Restore the default (Normal) style so no stray number-format style
# lingers on these cells once the literal text has been written.
$editRange.Style = "Normal"
